# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-rank two pairs of countries whose cumulative totals crossed over ---
# Costa Rica overtakes Portugal (row 50/51)
$ws.Range("A50").Value = "Costa Rica"
$ws.Range("A51").Value = "Portugal"

# Angola overtakes Guadalupe (row 119/120)
$ws.Range("A119").Value = "Angola"
$ws.Range("A120").Value = "Guadalupe"

# --- Refresh the "last updated" footer text ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 21:59"

# --- Refresh per-country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 8024612
$ws.Range("C4").Value = 32614
$ws.Range("D4").Value = 5163662
$ws.Range("E4").Value = 2641001
$ws.Range("G4").Value = 254
$ws.Range("H4").Value = 219949

$ws.Range("B5").Value = 7173345
$ws.Range("C5").Value = 54045
$ws.Range("D5").Value = 6224621
$ws.Range("E5").Value = 838830
$ws.Range("G5").Value = 710
$ws.Range("H5").Value = 109894

$ws.Range("B14").Value = 693359
$ws.Range("C14").Value = 888
$ws.Range("D14").Value = 624659
$ws.Range("E14").Value = 50837
$ws.Range("G14").Value = 83
$ws.Range("H14").Value = 17863

$ws.Range("B25").Value = 331057
$ws.Range("C25").Value = 4766
$ws.Range("E25").Value = 44436
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 9721

$ws.Range("B27").Value = 294031
$ws.Range("C27").Value = 3538
$ws.Range("D27").Value = 238681
$ws.Range("E27").Value = 53329
$ws.Range("G27").Value = 41
$ws.Range("H27").Value = 2021

$ws.Range("B34").Value = 147315
$ws.Range("C34").Value = 282
$ws.Range("E34").Value = 6963
$ws.Range("G34").Value = 27
$ws.Range("H34").Value = 12218

$ws.Range("B50").Value = 89223
$ws.Range("C50").Value = 733
$ws.Range("D50").Value = 53247
$ws.Range("E50").Value = 34868
$ws.Range("G50").Value = 19
$ws.Range("H50").Value = 1108

$ws.Range("B51").Value = 87913
$ws.Range("C51").Value = 1249
$ws.Range("D51").Value = 53498
$ws.Range("E51").Value = 32321
$ws.Range("G51").Value = 14
$ws.Range("H51").Value = 2094

$ws.Range("B66").Value = 53225
$ws.Range("C66").Value = 153
$ws.Range("E66").Value = 14034

$ws.Range("B92").Value = 20155
$ws.Range("C92").Value = 1
$ws.Range("D92").Value = 19831
$ws.Range("E92").Value = 204

$ws.Range("B119").Value = 6488
$ws.Range("C119").Value = 122
$ws.Range("D119").Value = 2744
$ws.Range("E119").Value = 3525
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 219

$ws.Range("B120").Value = 6483
$ws.Range("D120").Value = 2199
$ws.Range("E120").Value = 4207
$ws.Range("H120").Value = 77

$ws.Range("B123").Value = 5824
$ws.Range("C123").Value = 3
$ws.Range("D123").Value = 4659
$ws.Range("E123").Value = 985

$ws.Range("B133").Value = 4905
$ws.Range("C133").Value = 9
$ws.Range("D133").Value = 3877
$ws.Range("E133").Value = 996
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 32

$ws.Range("B148").Value = 3296
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 2533
$ws.Range("E148").Value = 631

$ws.Range("E159").Value = 127
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 596

$ws.Range("B165").Value = 1371
$ws.Range("C165").Value = 8
$ws.Range("D165").Value = 1250
$ws.Range("E165").Value = 39

$ws.Range("B166").Value = 1308
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 1117
$ws.Range("E166").Value = 99

$ws.Range("D181").Value = 463
$ws.Range("E181").Value = 14
